# Clear column A (rows 3:39) on the "Main" sheet - the per-module numeric/
# status markers that used to live in column A are removed entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("A3:A39").ClearContents()

# Update the active selection to match the saved view state.
$ws.Range("D43").Select()
